# Feria Lagunitas de Puerto Montt - Cilantro
# A new weekly price observation is inserted as a new data row (row 65),
# pushing every subsequent row down by one (old row 65 -> new row 66, ...,
# old row 146 -> new row 147). The worksheet dimension grows from
# A1:R146 to A1:R147 automatically once the new row is populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 65; this shifts rows 65..146 down to 66..147
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new observation.
# Columns that are constant for every record in this sheet (A, B, C, E, F,
# G, H, I, R) are simply repeated.
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 44413
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = 100112040
$ws.Range("G65").Value = "Cilantro"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 140
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 16000
$ws.Range("M65").Value = 15500
$ws.Range("N65").Value = "$/caja 36 atados"
$ws.Range("O65").Value = "Región Metropolitana"
$ws.Range("P65").Value = 431
$ws.Range("Q65").Value = 36
$ws.Range("R65").Value = "Hortaliza"
